{"js": "// Update the generated three-digit x one-digit multiplication practice\n// answers to the new set of problems/answers.\nconst replacements = [\n  [\"212\u00d74=848\", \"299\u00d77=2093\"],\n  [\"680\u00d78=5440\", \"531\u00d74=2124\"],\n  [\"733\u00d78=5864\", \"765\u00d77=5355\"],\n  [\"680\u00d77=4760\", \"759\u00d79=6831\"],\n  [\"853\u00d74=3412\", \"531\u00d76=3186\"],\n  [\"704\u00d72=1408\", \"459\u00d76=2754\"],\n  [\"944\u00d79=8496\", \"683\u00d73=2049\"],\n  [\"323\u00d79=2907\", \"687\u00d74=2748\"],\n  [\"500\u00d75=2500\", \"826\u00d75=4130\"],\n  [\"163\u00d78=1304\", \"978\u00d72=1956\"],\n  [\"226\u00d73=678\", \"827\u00d73=2481\"],\n  [\"656\u00d74=2624\", \"933\u00d77=6531\"],\n  [\"252\u00d76=1512\", \"402\u00d77=2814\"],\n  [\"664\u00d72=1328\", \"357\u00d77=2499\"],\n  [\"899\u00d79=8091\", \"285\u00d74=1140\"],\n  [\"508\u00d78=4064\", \"884\u00d74=3536\"],\n  [\"180\u00d76=1080\", \"693\u00d73=2079\"],\n  [\"962\u00d72=1924\", \"946\u00d78=7568\"],\n  [\"959\u00d78=7672\", \"344\u00d75=1720\"],\n  [\"702\u00d75=3510\", \"560\u00d74=2240\"],\n  [\"236\u00d79=2124\", \"601\u00d78=4808\"],\n  [\"640\u00d74=2560\", \"901\u00d73=2703\"],\n  [\"349\u00d77=2443\", \"494\u00d77=3458\"],\n  [\"939\u00d73=2817\", \"921\u00d74=3684\"],\n  [\"777\u00d78=6216\", \"150\u00d75=750\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const range of results.items) {\n    range.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Update the generated three-digit x one-digit multiplication practice\n# answers to the new set of problems/answers.\n$d = $word.ActiveDocument\n\n$replacements = @{\n    \"212\u00d74=848\"   = \"299\u00d77=2093\"\n    \"680\u00d78=5440\"  = \"531\u00d74=2124\"\n    \"733\u00d78=5864\"  = \"765\u00d77=5355\"\n    \"680\u00d77=4760\"  = \"759\u00d79=6831\"\n    \"853\u00d74=3412\"  = \"531\u00d76=3186\"\n    \"704\u00d72=1408\"  = \"459\u00d76=2754\"\n    \"944\u00d79=8496\"  = \"683\u00d73=2049\"\n    \"323\u00d79=2907\"  = \"687\u00d74=2748\"\n    \"500\u00d75=2500\"  = \"826\u00d75=4130\"\n    \"163\u00d78=1304\"  = \"978\u00d72=1956\"\n    \"226\u00d73=678\"   = \"827\u00d73=2481\"\n    \"656\u00d74=2624\"  = \"933\u00d77=6531\"\n    \"252\u00d76=1512\"  = \"402\u00d77=2814\"\n    \"664\u00d72=1328\"  = \"357\u00d77=2499\"\n    \"899\u00d79=8091\"  = \"285\u00d74=1140\"\n    \"508\u00d78=4064\"  = \"884\u00d74=3536\"\n    \"180\u00d76=1080\"  = \"693\u00d73=2079\"\n    \"962\u00d72=1924\"  = \"946\u00d78=7568\"\n    \"959\u00d78=7672\"  = \"344\u00d75=1720\"\n    \"702\u00d75=3510\"  = \"560\u00d74=2240\"\n    \"236\u00d79=2124\"  = \"601\u00d78=4808\"\n    \"640\u00d74=2560\"  = \"901\u00d73=2703\"\n    \"349\u00d77=2443\"  = \"494\u00d77=3458\"\n    \"939\u00d73=2817\"  = \"921\u00d74=3684\"\n    \"777\u00d78=6216\"  = \"150\u00d75=750\"\n}\n\nforeach ($old in $replacements.Keys) {\n    $new = $replacements[$old]\n    $rng = $d.Content\n    $find = $rng.Find\n    $find.ClearFormatting()\n    $find.Text = $old\n    $find.Replacement.Text = $new\n    $find.Forward = $true\n    $find.Wrap = 0\n    $find.MatchCase = $true\n    $find.MatchWholeWord = $false\n    $find.MatchWildcards = $false\n    $find.Execute($find.Text, $find.MatchCase, $find.MatchWholeWord, $find.MatchWildcards, $false, $false, $find.Forward, $find.Wrap, $false, $find.Replacement.Text, 2)\n}\n"}
